$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data in row 3 and row 4 represent the same market entry but for two
# different weekly dates. The edit swaps which date (and its associated
# volume/price figures) sits in row 3 vs row 4 ("Fruta / hortaliza, semanal").

# Capture current (pre-swap) values for the columns that differ between the
# two rows: D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado), P (Precio $/Kg).
$cols = @("D", "J", "K", "L", "M", "P")

$row3 = @{}
$row4 = @{}
foreach ($col in $cols) {
    $row3[$col] = $ws.Range("${col}3").Value2
    $row4[$col] = $ws.Range("${col}4").Value2
}

# Write row 4's original values into row 3, and row 3's original values into row 4.
foreach ($col in $cols) {
    $ws.Range("${col}3").Value = $row4[$col]
    $ws.Range("${col}4").Value = $row3[$col]
}
